$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be auto-detected as numbers
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "51.471.15"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.977.56"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "380.49"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "104.67"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").Value = "37.23"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "0.0845"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "3.448.07"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "18.54"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "7.49"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "2.973.28"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "0.975"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "51.450.99"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "12.99"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "69.10"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "262.77"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  +3.69%  "
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "7.19"
$ws.Range("E27").Value = "  +16.55%  "
$ws.Range("D28").Value = "7.46"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.113"
$ws.Range("E30").Value = "  +8.41%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "26.02"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "34.90"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "2.09"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "51.19"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").Value = "0.0455"
$ws.Range("E36").Value = "  +6.97%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "17.44"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("D40").Value = "2.59"
$ws.Range("E40").Value = "  -5.68%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").Value = "123.91"
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("D44").Value = "22.21"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "0.281"
$ws.Range("E45").Value = "  +18.84%  "
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("D48").Value = "2.038.00"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "0.0356"
$ws.Range("E50").Value = "  +11.83%  "
$ws.Range("D51").Value = "5.17"
$ws.Range("E51").Value = "  +2.76%  "
